$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.935.18'
$ws.Range('E2').Value = '  -2.58%  '
$ws.Range('D3').Value = '2.989.87'
$ws.Range('E3').Value = '  -1.83%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.33%  '
$ws.Range('D5').Value = '''528.64'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.06%  '
$ws.Range('D6').Value = '''131.58'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.58%  '
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('D8').Value = '2.986.94'
$ws.Range('E8').Value = '  -1.48%  '
$ws.Range('D9').Value = '''0.490'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -0.95%  '
$ws.Range('D10').Value = '''0.149'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -2.95%  '
$ws.Range('D11').Value = '''6.11'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -0.73%  '
$ws.Range('D12').Value = '''0.441'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -2.15%  '
$ws.Range('D13').Value = '''0.0000219'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -1.96%  '
$ws.Range('D14').Value = '''33.56'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.53%  '
$ws.Range('D15').Value = '3.468.87'
$ws.Range('E15').Value = '  -2.20%  '
$ws.Range('E16').Value = '  -0.01%  '
$ws.Range('D17').Value = '60.889.75'
$ws.Range('E17').Value = '  -2.77%  '
$ws.Range('D18').Value = '2.987.03'
$ws.Range('E18').Value = '  -2.29%  '
$ws.Range('D19').Value = '''6.52'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.85%  '
$ws.Range('D20').Value = '''462.14'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -3.86%  '
$ws.Range('D21').Value = '''13.12'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.47%  '
$ws.Range('D22').Value = '''0.669'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -3.01%  '
$ws.Range('D23').Value = '''6.85'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -3.34%  '
$ws.Range('D24').Value = '''78.76'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -0.21%  '
$ws.Range('D25').Value = '''11.83'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -2.07%  '
$ws.Range('E26').Value = '  +0.15%  '
$ws.Range('D27').Value = '''2.65'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.38%  '
$ws.Range('D28').Value = '''7.71'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -4.18%  '
$ws.Range('D29').Value = '''0.998'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.30%  '
$ws.Range('E30').Value = '  +3.26%  '
$ws.Range('D31').Value = '''1.86'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.19%  '
$ws.Range('D32').Value = '''25.29'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -1.85%  '
$ws.Range('D33').Value = '''55.04'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -2.28%  '
$ws.Range('E34').Value = '  -4.50%  '
$ws.Range('D35').Value = '''5.36'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('D36').Value = '''5.79'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -1.95%  '
$ws.Range('D37').Value = '''456.63'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -3.84%  '
$ws.Range('D38').Value = '3.175.17'
$ws.Range('E38').Value = '  +3.24%  '
$ws.Range('D39').Value = '''0.0780'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.51%  '
$ws.Range('D40').Value = '''0.0381'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -3.47%  '
$ws.Range('E41').Value = '  +1.34%  '
$ws.Range('D42').Value = '''8.06'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.17%  '
$ws.Range('D43').Value = '''2.43'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -7.75%  '
$ws.Range('D45').Value = '''26.04'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +6.26%  '
$ws.Range('D46').Value = '''0.244'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -3.07%  '
$ws.Range('D47').Value = '''118.54'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.55%  '
$ws.Range('D48').Value = '''0.108'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +0.21%  '
$ws.Range('E49').Value = '  -2.67%  '
$ws.Range('D50').Value = '0.0₃0492'
$ws.Range('E50').Value = '  -9.38%  '
$ws.Range('E51').Value = '  +7.74%  '
